# Auto-generated script applying the cryptos.xlsx price/volume update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cells whose new text would otherwise be auto-parsed as a number by Excel ---
# Mark them as Text format first so the literal string is preserved exactly.
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# --- Row-by-row Price (D) / Volume(1h) (E) updates ---
$ws.Range("D2").Value = "64.239.32"
$ws.Range("E2").Value = "  +2.82%  "
$ws.Range("D3").Value = "3.114.03"
$ws.Range("E3").Value = "  +2.77%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "546.09"
$ws.Range("E5").Value = "  +1.49%  "
$ws.Range("D6").Value = "141.50"
$ws.Range("E6").Value = "  +6.79%  "
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("D8").Value = "3.110.00"
$ws.Range("E8").Value = "  +2.75%  "
$ws.Range("E9").Value = "  +3.38%  "
$ws.Range("D10").Value = "6.58"
$ws.Range("E10").Value = "  +3.33%  "
$ws.Range("E11").Value = "  +3.64%  "
$ws.Range("E12").Value = "  +2.25%  "
$ws.Range("E13").Value = "  +8.04%  "
$ws.Range("D14").Value = "35.23"
$ws.Range("E14").Value = "  +3.11%  "
$ws.Range("D15").Value = "3.624.00"
$ws.Range("E15").Value = "  +3.04%  "
$ws.Range("D16").Value = "64.231.38"
$ws.Range("E16").Value = "  +2.73%  "
$ws.Range("E17").Value = "  +2.61%  "
$ws.Range("D18").Value = "3.114.94"
$ws.Range("E18").Value = "  +3.12%  "
$ws.Range("D19").Value = "6.73"
$ws.Range("E19").Value = "  +3.21%  "
$ws.Range("D20").Value = "488.40"
$ws.Range("E20").Value = "  +2.61%  "
$ws.Range("D21").Value = "13.55"
$ws.Range("E21").Value = "  +2.53%  "
$ws.Range("D22").Value = "0.709"
$ws.Range("E22").Value = "  +3.14%  "
$ws.Range("D23").Value = "7.22"
$ws.Range("E23").Value = "  +3.68%  "
$ws.Range("D24").Value = "79.86"
$ws.Range("E24").Value = "  +4.61%  "
$ws.Range("D25").Value = "12.42"
$ws.Range("E25").Value = "  +3.48%  "
$ws.Range("E26").Value = "  +0.32%  "
$ws.Range("D27").Value = "2.76"
$ws.Range("E27").Value = "  +3.97%  "
$ws.Range("D28").Value = "8.27"
$ws.Range("E28").Value = "  +2.41%  "
$ws.Range("D29").Value = "0.998"
$ws.Range("E29").Value = "  -0.11%  "
$ws.Range("D30").Value = "26.59"
$ws.Range("E30").Value = "  +2.66%  "
$ws.Range("D31").Value = "1.93"
$ws.Range("E31").Value = "  +1.66%  "
$ws.Range("D32").Value = "1.16"
$ws.Range("E32").Value = "  +3.92%  "
$ws.Range("D33").Value = "2.40"
$ws.Range("E33").Value = "  -1.84%  "
$ws.Range("D34").Value = "57.78"
$ws.Range("E34").Value = "  -3.46%  "
$ws.Range("D35").Value = "505.30"
$ws.Range("E35").Value = "  -0.63%  "
$ws.Range("E36").Value = "  +7.70%  "
$ws.Range("E37").Value = "  +4.37%  "
$ws.Range("D38").Value = "3.288.49"
$ws.Range("E38").Value = "  +8.79%  "
$ws.Range("D39").Value = "0.0406"
$ws.Range("E39").Value = "  +3.55%  "
$ws.Range("D40").Value = "0.0805"
$ws.Range("E40").Value = "  +3.80%  "
$ws.Range("D41").Value = "0.122"
$ws.Range("E41").Value = "  +4.75%  "
$ws.Range("D42").Value = "2.75"
$ws.Range("E42").Value = "  +8.04%  "
$ws.Range("D43").Value = "8.18"
$ws.Range("E43").Value = "  +3.17%  "
$ws.Range("D44").Value = "0.259"
$ws.Range("E44").Value = "  +4.71%  "
$ws.Range("E45").Value = "  +0.08%  "
$ws.Range("D46").Value = "2.08"
$ws.Range("E46").Value = "  +4.69%  "
$ws.Range("D49").Value = "25.21"
$ws.Range("E49").Value = "  +5.73%  "
$ws.Range("E50").Value = "  +4.12%  "
$ws.Range("D51").Value = "2.45"
$ws.Range("E51").Value = "  +4.59%  "

# --- Rows 47 & 48 swap places (PEPE <-> Monero) with refreshed data ---
$ws.Range("B47").Value = "Monero"
$ws.Range("C47").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D47").Value = "123.72"
$ws.Range("E47").Value = "  +4.67%  "
$ws.Range("B48").Value = "PEPE"
$ws.Range("C48").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D48").Value = "0.0₃0541"
$ws.Range("E48").Value = "  +12.43%  "
